$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.026.24"
$ws.Range("E2").Value = "  -4.80%  "

$ws.Range("D3").Value = "3.275.69"
$ws.Range("E3").Value = "  -5.63%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.26"
$ws.Range("E5").Value = "  -3.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.20"
$ws.Range("E6").Value = "  -3.11%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("E8").Value = "  -2.64%  "

$ws.Range("D9").Value = "3.268.10"
$ws.Range("E9").Value = "  -5.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.188"
$ws.Range("E10").Value = "  -8.18%  "

$ws.Range("E11").Value = "  -4.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.50"
$ws.Range("E12").Value = "  -7.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000267"
$ws.Range("E13").Value = "  -6.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.63"
$ws.Range("E14").Value = "  -5.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "631.91"
$ws.Range("E15").Value = "  -1.37%  "

$ws.Range("D16").Value = "3.799.25"
$ws.Range("E16").Value = "  -5.72%  "

$ws.Range("D17").Value = "65.902.46"
$ws.Range("E17").Value = "  -4.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.90"
$ws.Range("E18").Value = "  -1.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.117"
$ws.Range("E19").Value = "  -3.43%  "

$ws.Range("D20").Value = "3.277.09"
$ws.Range("E20").Value = "  -5.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.36"
$ws.Range("E21").Value = "  -7.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.904"
$ws.Range("E22").Value = "  -3.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.35"
$ws.Range("E23").Value = "  +2.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "106.84"
$ws.Range("E24").Value = "  +7.94%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.91"
$ws.Range("E25").Value = "  -7.14%  "

$ws.Range("E26").Value = "  -7.20%  "

$ws.Range("E27").Value = "  -6.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.61"
$ws.Range("E28").Value = "  -2.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.71"
$ws.Range("E29").Value = "  -6.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.33"
$ws.Range("E30").Value = "  -6.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.05"
$ws.Range("E31").Value = "  -5.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.26"
$ws.Range("E32").Value = "  -6.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.05"
$ws.Range("E33").Value = "  -4.67%  "

$ws.Range("E34").Value = "  -3.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "532.53"
$ws.Range("E35").Value = "  +3.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.51"
$ws.Range("E36").Value = "  -5.57%  "

$ws.Range("D37").Value = "3.727.28"
$ws.Range("E37").Value = "  +1.11%  "

$ws.Range("E38").Value = "  +0.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("E39").Value = "  -2.88%  "

$ws.Range("D40").Value = "0.0₃0731"
$ws.Range("E40").Value = "  -7.38%  "

$ws.Range("E41").Value = "  -1.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.75"
$ws.Range("E42").Value = "  -6.27%  "

$ws.Range("E43").Value = "  -1.97%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "32.71"
$ws.Range("E44").Value = "  -4.71%  "

$ws.Range("E45").Value = "  -8.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.29"
$ws.Range("E46").Value = "  -1.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0415"
$ws.Range("E47").Value = "  -6.34%  "

$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.61"
$ws.Range("E48").Value = "  -7.04%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.129"
$ws.Range("E49").Value = "  -3.36%  "

$ws.Range("E50").Value = "  -0.16%  "

$ws.Range("E51").Value = "  +1.87%  "
